$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numeric-looking strings (prices, percentages) that must
# stay TEXT, not be auto-coerced to numbers by Excel. Force text format,
# assign, then clear the formatting residue so cell styling matches source.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '70.095.07'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.565.04'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.72%  '
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '198.65'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +1.38%  '
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '587.68'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -2.73%  '
$c.ClearFormats()
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.622'
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.213'
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +3.96%  '
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.638'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -2.15%  '
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '52.81'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.00%  '
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000298'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -2.21%  '
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '9.48'
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '692.98'
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +15.66%  '
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.141.11'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '70.206.05'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '12.61'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -3.41%  '
$c.ClearFormats()
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '18.86'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -2.16%  '
$c.ClearFormats()
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.564.13'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.69%  '
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.984'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '18.17'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +1.59%  '
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '110.54'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +7.43%  '
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.26'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +2.29%  '
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '4.47'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -3.67%  '
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -3.66%  '
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.42'
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -3.94%  '
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.01'
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.87'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +2.87%  '
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '34.15'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.42'
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +3.38%  '
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.01'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -1.51%  '
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '12.11'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -2.15%  '
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -2.58%  '
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '63.31'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0₃0836'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +1.63%  '
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.799.94'
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -2.44%  '
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +0.99%  '
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '504.90'
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -4.27%  '
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -8.33%  '
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '36.00'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.378'
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -3.83%  '
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +2.12%  '
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0466'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +2.47%  '
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +5.62%  '
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.43'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +3.31%  '
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.36%  '
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '8.56'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.38%  '
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.32%  '
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.79'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +20.94%  '
$c.ClearFormats()
